# Weekly update: insert a new daily price record for
# "Feria Lagunitas de Puerto Montt - Arveja Verde" at row 32,
# pushing the existing rows 32:98 down to 33:99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 32 (shifts rows 32-98 down to 33-99,
# Excel also grows the sheet dimension from A1:R98 to A1:R99 automatically).
$ws.Rows(32).Insert()

# Populate the newly inserted row 32 with the new record's data.
$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C32").Value = "Los Lagos"
$ws.Range("D32").Value = 44544
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 100112022
$ws.Range("G32").Value = "Arveja Verde"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 20000
$ws.Range("N32").Value = '$/saco 25 kilos'
$ws.Range("O32").Value = "Región de La Araucanía"
$ws.Range("P32").Value = 800
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
